# The deck ships two theme parts:
#   - the "Integral" colour theme, which is the theme actually applied to the
#     slide master (and therefore to every slide in the deck)
#   - the stock "Office Theme" colour theme
# The commit swaps which colours are applied to the presentation's live
# theme: the slide master's theme colours change from the Integral palette
# to the default Office palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# Font scheme and format scheme (fills/lines/effects) are identical between
# the two theme parts, so only the twelve theme colours need to change.

$p = $ppt.ActivePresentation

# The slide master's Theme.ThemeColorScheme exposes the twelve colours of
# the theme that is actually rendered behind every slide, in the fixed
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$sm  = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# Target palette = the stock "Office Theme" colours (RGB hex -> value).
$officeTheme = @{
    1  = 0x000000  # dk1
    2  = 0xFFFFFF  # lt1
    3  = 0x44546A  # dk2
    4  = 0xE7E6E6  # lt2
    5  = 0x5B9BD5  # accent1
    6  = 0xED7D31  # accent2
    7  = 0xA5A5A5  # accent3
    8  = 0xFFC000  # accent4
    9  = 0x4472C4  # accent5
    10 = 0x70AD47  # accent6
    11 = 0x0563C1  # hlink
    12 = 0x954F72  # folHlink
}

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeTheme[$i]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    # PowerPoint's RGB colour values are packed as 0x00BBGGRR.
    $tcs.Item($i).RGB = ($b * 0x10000) + ($g * 0x100) + $r
}
